$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) values
$ws.Range("D2").Value = "26.704.14"
$ws.Range("D3").Value = "1.598.04"
$ws.Range("D5").Value = "'211.57"
$ws.Range("D10").Value = "'19.48"
$ws.Range("D12").Value = "1.823.61"
$ws.Range("D13").Value = "1.591.20"
$ws.Range("D16").Value = "'65.19"
$ws.Range("D17").Value = "26.665.67"
$ws.Range("D19").Value = "'209.53"
$ws.Range("D21").Value = "'7.08"
$ws.Range("D24").Value = "'8.94"
$ws.Range("D25").Value = "'143.15"
$ws.Range("D29").Value = "'15.31"
$ws.Range("D30").Value = "'0.0519"
$ws.Range("D31").Value = "'1.15"
$ws.Range("D34").Value = "1.285.51"
$ws.Range("D35").Value = "'0.618"
$ws.Range("D41").Value = "'5.44"
$ws.Range("D43").Value = "'2.18"
$ws.Range("D44").Value = "'63.14"
$ws.Range("D45").Value = "1.735.77"
$ws.Range("D46").Value = "'91.20"
$ws.Range("D49").Value = "'0.0507"

# Update Volume(1h) (column E) values
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("E9").Value = "  +1.07%  "
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("E15").Value = "  +0.61%  "
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("E18").Value = "  +5.06%  "
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  +4.71%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("E25").Value = "  -1.67%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +0.62%  "
$ws.Range("E30").Value = "  +3.20%  "
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  -6.89%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  +16.93%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("E46").Value = "  +1.58%  "
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("E51").Value = "  -1.15%  "

# Reset style on Price cells that received a quote-prefix so the cell style stays default
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D49").Style = "Normal"
